$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 272, shifting existing rows 272-326 down to 273-327
$ws.Rows.Item(272).Insert()

# Populate the newly inserted row 272 with the new record's data
$ws.Cells.Item(272, 1).Value = 10
$ws.Cells.Item(272, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(272, 3).Value = "La Araucanía"
$ws.Cells.Item(272, 4).Value = 45209
$ws.Cells.Item(272, 5).Value = 9
$ws.Cells.Item(272, 6).Value = "Fruta"
$ws.Cells.Item(272, 7).Value = 100101
$ws.Cells.Item(272, 8).Value = "Berries"
$ws.Cells.Item(272, 9).Value = 100112025
$ws.Cells.Item(272, 10).Value = "Frutilla"
$ws.Cells.Item(272, 11).Value = "Sin especificar"
$ws.Cells.Item(272, 12).Value = "Primera"
$ws.Cells.Item(272, 13).Value = 35
$ws.Cells.Item(272, 14).Value = 14000
$ws.Cells.Item(272, 15).Value = 14000
$ws.Cells.Item(272, 16).Value = 14000
$ws.Cells.Item(272, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(272, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(272, 19).Value = 2000
$ws.Cells.Item(272, 20).Value = 7

# Match the date number format used by the rest of column D (style index 2 -> YYYY-MM-DD HH:MM:SS)
$ws.Cells.Item(272, 4).NumberFormat = $ws.Cells.Item(273, 4).NumberFormat
